$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 609.2
$ws.Range("I33").Value = 269.9091
$ws.Range("J33").Value = 875.7857
$ws.Range("K33").Value = 269.9091
$ws.Range("L33").Value = 875.7857
$ws.Range("M33").Value = -40.90910000000002
$ws.Range("N33").Value = -1333.7857
$ws.Range("H51").Value = 15108.952
$ws.Range("I51").Value = 23780
$ws.Range("J51").Value = 12399.25
$ws.Range("K51").Value = 23780
$ws.Range("L51").Value = 12399.25
$ws.Range("M51").Value = -23296
$ws.Range("N51").Value = -13367.25
$ws.Range("H53").Value = 990.11536
$ws.Range("I53").Value = 941.25
$ws.Range("J53").Value = 1068.3
$ws.Range("K53").Value = 941.25
$ws.Range("L53").Value = 1068.3
$ws.Range("M53").Value = -304.25
$ws.Range("N53").Value = -2342.3
$ws.Range("H69").Value = 9980
$ws.Range("I69").Value = 9976
$ws.Range("K69").Value = 29928
$ws.Range("M69").Value = -29054
$ws.Range("H72").Value = 9980
$ws.Range("I72").Value = 9976
$ws.Range("K72").Value = 89784
$ws.Range("M72").Value = -85416
$ws.Range("H86").Value = 41656.96
$ws.Range("I86").Value = 46268.78
$ws.Range("K86").Value = 46268.78
$ws.Range("M86").Value = -45145.78
$ws.Range("H89").Value = 41656.96
$ws.Range("I89").Value = 46268.78
$ws.Range("K89").Value = 231343.9
$ws.Range("M89").Value = -225727.9
$ws.Range("H132").Value = 1007.5476
$ws.Range("I132").Value = 961.1625
$ws.Range("K132").Value = 2883.4875
$ws.Range("M132").Value = -353.4875000000002
$ws.Range("H135").Value = 1166.5814
$ws.Range("I135").Value = 916.6111
$ws.Range("J135").Value = 2452.1428
$ws.Range("K135").Value = 8249.499899999999
$ws.Range("L135").Value = 22069.2852
$ws.Range("M135").Value = -5714.499899999999
$ws.Range("N135").Value = -27139.2852
$ws.Range("H138").Value = 2180.2173
$ws.Range("I138").Value = 898.9655
$ws.Range("J138").Value = 2770
$ws.Range("K138").Value = 2696.8965
$ws.Range("L138").Value = 8310
$ws.Range("M138").Value = 2443.1035
$ws.Range("N138").Value = -18590
$ws.Range("H141").Value = 1935.909
$ws.Range("I141").Value = 835.52
$ws.Range("J141").Value = 5374.625
$ws.Range("K141").Value = 2506.56
$ws.Range("L141").Value = 16123.875
$ws.Range("M141").Value = 2673.44
$ws.Range("N141").Value = -26483.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3140.1843
$ws.Range("I61").Value = 1851.1
$ws.Range("K61").Value = 1851.1
$ws.Range("M61").Value = -1639.1
$ws.Range("H74").Value = 145206
$ws.Range("I74").Value = 168158.17
$ws.Range("J74").Value = 7493
$ws.Range("K74").Value = 168158.17
$ws.Range("L74").Value = 7493
$ws.Range("M74").Value = -167284.17
$ws.Range("N74").Value = -9241
$ws.Range("H77").Value = 145206
$ws.Range("I77").Value = 168158.17
$ws.Range("J77").Value = 7493
$ws.Range("K77").Value = 840790.8500000001
$ws.Range("L77").Value = 37465
$ws.Range("M77").Value = -836422.8500000001
$ws.Range("N77").Value = -46201
$ws.Range("H88").Value = 1757.875
$ws.Range("I88").Value = 1841.5834
$ws.Range("J88").Value = 1506.75
$ws.Range("K88").Value = 1841.5834
$ws.Range("L88").Value = 1506.75
$ws.Range("M88").Value = -1435.5834
$ws.Range("N88").Value = -2318.75
$ws.Range("H91").Value = 1757.875
$ws.Range("I91").Value = 1841.5834
$ws.Range("J91").Value = 1506.75
$ws.Range("K91").Value = 1841.5834
$ws.Range("L91").Value = 1506.75
$ws.Range("M91").Value = -437.5834
$ws.Range("N91").Value = -4314.75
$ws.Range("H97").Value = 2558.6
$ws.Range("I97").Value = 2129.1667
$ws.Range("J97").Value = 3202.75
$ws.Range("K97").Value = 2129.1667
$ws.Range("L97").Value = 3202.75
$ws.Range("M97").Value = -1633.1667
$ws.Range("N97").Value = -4194.75
$ws.Range("H136").Value = 3140.1843
$ws.Range("I136").Value = 1851.1
$ws.Range("K136").Value = 5553.299999999999
$ws.Range("M136").Value = -3003.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 34995
$ws.Range("J40").Value = 34995
$ws.Range("L40").Value = 34995
$ws.Range("N40").Value = -35525
$ws.Range("H105").Value = 2776.4666
$ws.Range("I105").Value = 2800.4285
$ws.Range("J105").Value = 2441
$ws.Range("K105").Value = 2800.4285
$ws.Range("L105").Value = 2441
$ws.Range("M105").Value = -1053.4285
$ws.Range("N105").Value = -5935
$ws.Range("H134").Value = 1870.1212
$ws.Range("I134").Value = 1913.25
$ws.Range("J134").Value = 490
$ws.Range("K134").Value = 5739.75
$ws.Range("L134").Value = 1470
$ws.Range("M134").Value = -3204.75
$ws.Range("N134").Value = -6540

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 219568.73
$ws.Range("J31").Value = 2587.3333
$ws.Range("L31").Value = 2587.3333
$ws.Range("N31").Value = -3177.3333
$ws.Range("H34").Value = 219568.73
$ws.Range("J34").Value = 2587.3333
$ws.Range("L34").Value = 2587.3333
$ws.Range("N34").Value = -2991.3333
$ws.Range("H58").Value = 1772.0769
$ws.Range("I58").Value = 1772.0769
$ws.Range("K58").Value = 1772.0769
$ws.Range("M58").Value = -1569.0769
$ws.Range("H105").Value = 5354.467
$ws.Range("I105").Value = 2028.4166
$ws.Range("J105").Value = 7571.8335
$ws.Range("K105").Value = 2028.4166
$ws.Range("L105").Value = 7571.8335
$ws.Range("M105").Value = -281.4166
$ws.Range("N105").Value = -11065.8335
$ws.Range("H132").Value = 1893.6154
$ws.Range("I132").Value = 1893.6154
$ws.Range("K132").Value = 5680.8462
$ws.Range("M132").Value = -3150.8462
$ws.Range("H134").Value = 2556.4644
$ws.Range("I134").Value = 2568.4363
$ws.Range("J134").Value = 1898
$ws.Range("K134").Value = 7705.3089
$ws.Range("L134").Value = 5694
$ws.Range("M134").Value = -5170.3089
$ws.Range("N134").Value = -10764
$ws.Range("H136").Value = 1772.0769
$ws.Range("I136").Value = 1772.0769
$ws.Range("K136").Value = 5316.2307
$ws.Range("M136").Value = -2766.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 392.19232
$ws.Range("I12").Value = 513.8889
$ws.Range("J12").Value = 327.7647
$ws.Range("K12").Value = 1541.6667
$ws.Range("L12").Value = 983.2941000000001
$ws.Range("M12").Value = -1368.6667
$ws.Range("N12").Value = -1329.2941
$ws.Range("H55").Value = 5868.1177
$ws.Range("I55").Value = 4327.7144
$ws.Range("K55").Value = 12983.1432
$ws.Range("M55").Value = -12806.1432
$ws.Range("H57").Value = 915
$ws.Range("I57").Value = 1940
$ws.Range("K57").Value = 5820
$ws.Range("M57").Value = -5261
$ws.Range("H68").Value = 8335684.5
$ws.Range("J68").Value = 3450.8333
$ws.Range("L68").Value = 10352.4999
$ws.Range("N68").Value = -11974.4999
$ws.Range("H71").Value = 8335684.5
$ws.Range("J71").Value = 3450.8333
$ws.Range("L71").Value = 31057.4997
$ws.Range("N71").Value = -39169.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7715.0586
$ws.Range("I70").Value = 8046.615
$ws.Range("J70").Value = 6637.5
$ws.Range("K70").Value = 8046.615
$ws.Range("L70").Value = 6637.5
$ws.Range("M70").Value = -7776.615
$ws.Range("N70").Value = -7177.5
$ws.Range("H73").Value = 7715.0586
$ws.Range("I73").Value = 8046.615
$ws.Range("J73").Value = 6637.5
$ws.Range("K73").Value = 8046.615
$ws.Range("L73").Value = 6637.5
$ws.Range("M73").Value = -7110.615
$ws.Range("N73").Value = -8509.5
$ws.Range("H122").Value = 1688.7778
$ws.Range("I122").Value = 1499.5
$ws.Range("K122").Value = 4498.5
$ws.Range("M122").Value = -2048.5
$ws.Range("H132").Value = 57622.54
$ws.Range("I132").Value = 65509
$ws.Range("J132").Value = 24499.4
$ws.Range("K132").Value = 196527
$ws.Range("L132").Value = 73498.20000000001
$ws.Range("M132").Value = -193997
$ws.Range("N132").Value = -78558.20000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 31714.143
$ws.Range("I29").Value = 28299.9
$ws.Range("K29").Value = 28299.9
$ws.Range("M29").Value = -28004.9
$ws.Range("H122").Value = 5538.909
$ws.Range("I122").Value = 4971.0713
$ws.Range("J122").Value = 6532.625
$ws.Range("K122").Value = 14913.2139
$ws.Range("L122").Value = 19597.875
$ws.Range("M122").Value = -12463.2139
$ws.Range("N122").Value = -24497.875
$ws.Range("H132").Value = 4754.8594
$ws.Range("I132").Value = 3357.3513
$ws.Range("J132").Value = 6669.963
$ws.Range("K132").Value = 10072.0539
$ws.Range("L132").Value = 20009.889
$ws.Range("M132").Value = -7542.053899999999
$ws.Range("N132").Value = -25069.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14715880
$ws.Range("J62").Value = 22738224
$ws.Range("L62").Value = 22738224
$ws.Range("N62").Value = -22739472
$ws.Range("H65").Value = 14715880
$ws.Range("J65").Value = 22738224
$ws.Range("L65").Value = 113691120
$ws.Range("N65").Value = -113697360
$ws.Range("H81").Value = 2208.923
$ws.Range("I81").Value = 1883.6818
$ws.Range("J81").Value = 3997.75
$ws.Range("K81").Value = 3767.3636
$ws.Range("L81").Value = 7995.5
$ws.Range("M81").Value = -2706.3636
$ws.Range("N81").Value = -10117.5
$ws.Range("H84").Value = 2208.923
$ws.Range("I84").Value = 1883.6818
$ws.Range("J84").Value = 3997.75
$ws.Range("K84").Value = 18836.818
$ws.Range("L84").Value = 39977.5
$ws.Range("M84").Value = -13532.818
$ws.Range("N84").Value = -50585.5
$ws.Range("H86").Value = 89974
$ws.Range("J86").Value = 89974
$ws.Range("L86").Value = 89974
$ws.Range("N86").Value = -92220
$ws.Range("H89").Value = 89974
$ws.Range("J89").Value = 89974
$ws.Range("L89").Value = 449870
$ws.Range("N89").Value = -461102
$ws.Range("H132").Value = 1775.1282
$ws.Range("I132").Value = 1648.5135
$ws.Range("K132").Value = 4945.5405
$ws.Range("M132").Value = -2415.5405
$ws.Range("H136").Value = 199698.64
$ws.Range("I136").Value = 235826.84
$ws.Range("J136").Value = 5509.625
$ws.Range("K136").Value = 707480.52
$ws.Range("L136").Value = 16528.875
$ws.Range("M136").Value = -704930.52
